# Rename the "ufs" sheet to "zero" and make it the active sheet/tab,
# matching the target workbook state:
#   - <sheet name="ufs" .../> -> <sheet name="zero" .../>
#   - "prot" sheet loses tabSelected, "zero" sheet gains tabSelected
#   - bookViews/workbookView gets activeTab="2" (zero-based index of "zero")

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ufs")
$ws.Name = "zero"

# Activating the renamed sheet updates the sheetView tabSelected flags
# (removes it from the previously active sheet, adds it here) and sets
# the workbook's bookViews/workbookView activeTab to this sheet's index.
$ws.Activate()

Write-Output "Renamed 'ufs' to 'zero' and activated it"
